$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last_edited_time" timestamp (shared across several rows in column D)
$newTimestamp = "2024-08-28T15:04:00.000Z"
foreach ($r in @(2, 3, 5, 6, 7, 8, 11, 13)) {
    $ws.Range("D$r").Value = $newTimestamp
}

# Update the "ứng lương" (salary advance) figures in row 7
$ws.Range("S7").Value = 67695000
$ws.Range("W7").Value = 121505000
$ws.Range("AE7").Value = 189200000
$ws.Range("AH7").Value = 172200000
$ws.Range("AK7").Value = 22
$ws.Range("AQ7").Value = 204200000
